# Apply symbol-list refresh (commit: "Updated symbol list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.34"
$ws.Range("D3").Value = "'23.22"
$ws.Range("D4").Value = "'6.472"
$ws.Range("D5").Value = "'0.06301"
$ws.Range("D6").Value = "'3.661"
$ws.Range("D7").Value = "'6.678"
$ws.Range("D8").Value = "'1.388"
$ws.Range("D9").Value = "'0.8355"
$ws.Range("D11").Value = "'0.1648"
$ws.Range("D12").Value = "'0.08299"
$ws.Range("D13").Value = "'0.03455"
$ws.Range("D14").Value = "'0.03110"
$ws.Range("D15").Value = "'0.09315"
$ws.Range("D16").Value = "'3.862"
$ws.Range("D17").Value = "'0.001657"
$ws.Range("D18").Value = "'0.04791"
$ws.Range("D19").Value = "'0.006312"
$ws.Range("D20").Value = "'0.005691"
$ws.Range("D21").Value = "'0.001084"
$ws.Range("D23").Value = "'3.712"
$ws.Range("D24").Value = "'2.370"
$ws.Range("D27").Value = "'0.0002680"
$ws.Range("D40").Value = "'0.04723"
$ws.Range("D41").Value = "'0.007059"
$ws.Range("D42").Value = "'0.1161"
$ws.Range("D43").Value = "'0.003699"
$ws.Range("D44").Value = "'0.01217"
$ws.Range("D45").Value = "'0.00006244"
$ws.Range("D48").Value = "'0.7966"
$ws.Range("B49").Value = "CryptobidCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("D49").Value = "'0.00002300"
$ws.Range("E49").Value = "48CryptobidCoinCBC"
$ws.Range("B50").Value = "BOLO"
$ws.Range("C50").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D50").Value = "'0.01389"
$ws.Range("E50").Value = "49BOLOBOLOBestin24h"
